$wb = $excel.ActiveWorkbook

$itSheet = $wb.Worksheets.Item("IT")
$aboutSheet = $wb.Worksheets.Item("About")

# Data update: advance/roll back the Initial Time value on the IT sheet.
$itSheet.Range("B2").Value = 2021

# Update cursor/selection state left in each sheet.
$itSheet.Range("B3").Select()
$aboutSheet.Range("A1").Select()

# "About" becomes the active (tab-selected) sheet.
$aboutSheet.Activate()
